# Add a "Total" column (F) = SUM(C:E) to the PBO and Service Cost sheets,
# matching the "updated data to include total plan" commit.

$wb = $excel.ActiveWorkbook

# --- PBO (sheet 1) and Service Cost (sheet 2): add Total column ---
foreach ($idx in 1, 2) {
    $ws = $wb.Worksheets.Item($idx)

    # Header
    $ws.Range("F1").Value = "Total"

    # First data row gets its own formula (becomes the "master" formula),
    # the rest are filled down from it so every row sums Retirement+Pension+IBT.
    $ws.Range("F2").Formula = "=SUM(C2:E2)"
    $ws.Range("F3:F81").FormulaR1C1 = "=SUM(RC[-3]:RC[-1])"
}

# Column F width on the PBO sheet (matches the other data columns' look)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(6).ColumnWidth = 10.1

# --- Restore / set the view state (active cell + selection) on each sheet ---

# PVFB (sheet 3): no data change, just cursor parked at H10
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("H10").Select()

# Service Cost (sheet 2): cursor parked at the bottom of the new column
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("F81").Select()

# PBO (sheet 1) ends up the active tab, with the new column selected
$ws1.Activate()
$ws1.Range("F1:F81").Select()
